$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 30.08. (introduce the new date string first so shared-string
# ordering matches the author's save)
$ws.Range("A14").Value = "30.08."

# 28.08. row's topic text is rewritten
$ws.Range("C12").Value = "Setup ubelix, Code: adjust reward, split in modules, add comments"

# 27.08. row's topic text is rewritten
$ws.Range("C11").Value = "Setup Ubelix, Code: REINFORCE algo, replay memory"

# 29.08. row now holds the "Meeting + Preparation" topic (moved up one slot)
$ws.Range("A13").Value = "29.08."
$ws.Range("B13").Value = 65
$ws.Range("C13").Value = "Meeting + Preparation"

# Finish populating the new 30.08. row
$ws.Range("B14").Value = 95
$ws.Range("C14").Value = "Code: Introduce num_message_passing, merge two gnn"

# Update selection to match the saved cursor position
$ws.Range("C15").Select()

$wb.Save()
